# Apply the commit's changes to the workbook:
# 1. Update the "Date" metadata value on the Metadata sheet.
# 2. Update Min / Max / Base Min / Base Max from "1" to "0" (stored as
#    text, like the rest of the column) for the ActorSystem.XCN9.composant1
#    row (row 7) on the Elements sheet.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B8").Value = "2025-05-05T11:54:16+00:00"

$wsElements = $wb.Worksheets.Item("Elements")
foreach ($addr in @("F7", "G7", "AG7", "AH7")) {
    $cell = $wsElements.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "0"
}
